# feat(stats): Actual statistical analysis
#
# - Removes the old two-sample t-test scratch area (COUNT/STDEV.S/AVERAGE/
#   SQRT/ABS/T.DIST helper columns K:R on rows 1-2 and 23-24, plus the
#   explanatory comment in N3/N25) from Sheet1.
# - Adds a new Sheet2 with a clean Year / Mean NDVI / Difference /
#   Relative Difference (normalized) table plus a small
#   AVERAGE / STDEV.S / ratio summary, replacing the removed analysis.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Drop the leftover "_xlchart.v1.*" hidden defined names.
# ---------------------------------------------------------------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# ---------------------------------------------------------------------
# 2. Clean up Sheet1's old stats scratch area (two blocks: rows 1-3 for
#    the additive-change group, rows 23-25 for the multiplicative-change
#    group). J1/J2/J23/J24 (labels "CHANGE"/"Additive Change"/
#    "Multiplicative Change") and the J3:J21 / J25:J43 change-series
#    stay untouched; only the K:R helper columns + the long comment in
#    column N go away.
# ---------------------------------------------------------------------
$ws1.Range("K1:R1").ClearContents()
$ws1.Range("K2:R2").ClearContents()
$ws1.Range("N3").ClearContents()

$ws1.Range("K23:R23").ClearContents()
$ws1.Range("K24:R24").ClearContents()
$ws1.Range("N25").ClearContents()

# Restore the selection Sheet1 should land on once it is no longer the
# active tab.
$ws1.Range("C2:C21").Select()

# ---------------------------------------------------------------------
# 3. Add Sheet2 (right after Sheet1) with the cleaned-up Year/NDVI table.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Year"
$ws2.Range("B1").Value = "Mean NDVI (0.0001 scale)"
$ws2.Range("C1").Value = "Difference"
$ws2.Range("D1").Value = "Relative Difference (normalized)"

# Copy Year (Sheet1 col H) / Mean NDVI (Sheet1 col C) for the 20 years,
# rows 2-21.
for ($i = 0; $i -lt 20; $i++) {
    $r = $i + 2
    $year = $ws1.Cells.Item($r, 8).Value2
    $ndvi = $ws1.Cells.Item($r, 3).Value2
    $ws2.Cells.Item($r, 1).Value = $year
    $ws2.Cells.Item($r, 2).Value = $ndvi
}

# Year-over-year absolute / relative differences.
$ws2.Range("C3").Formula = "=B3-B2"
$ws2.Range("C4:C21").Formula = "=B4-B3"
$ws2.Range("D3").Formula = "=C3/B2"
$ws2.Range("D4:D21").Formula = "=C4/B3"

# Small summary block.
$ws2.Range("F3").Formula = "=AVERAGE(C3:C21)"
$ws2.Range("F4").Formula = "=STDEV.S(C3:C21)"
$ws2.Range("G4").Formula = "=F3/F4"

# Number formats.
$ws2.Range("B2:B21").NumberFormat = "0"
$ws2.Range("C3:C21").NumberFormat = "0"
$ws2.Range("D3:D21").NumberFormat = "0%"
$ws2.Range("F3").NumberFormat = "0.00"

$ws2.Range("C1").ColumnWidth = 12.1640625

# Final selection on Sheet2 (also leaves Sheet2 as the active tab).
$ws2.Range("F3:G4").Select()
